# Updated C3DC phs000466 queries
#
# The "Treatment" tab query (cell B5) is rewritten so the treatment-agent
# column uses a plain REPLACE(...) instead of the redundant
# CONCAT(REPLACE(...)) wrapper. The other query cells (B4 Diagnosis,
# B6 TreatmentResp, B7 Survival) keep their existing text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$treatmentQuery = @"
SELECT
    DISTINCT prt.participant_id AS "Participant Id",
    trt.treatment_id AS "Treatment Id",
    CASE 
    WHEN trt.age_at_treatment_start = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_start >= 1000 THEN 
        substr(trt.age_at_treatment_start, 1, length(trt.age_at_treatment_start) - 3) || ',' || substr(trt.age_at_treatment_start, -3)
    ELSE 
        trt.age_at_treatment_start 
END AS "Age at Treatment Start",
    CASE 
    WHEN trt.age_at_treatment_end = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_end >= 1000 THEN 
        substr(trt.age_at_treatment_end, 1, length(trt.age_at_treatment_end) - 3) || ',' || substr(trt.age_at_treatment_end, -3)
    ELSE 
        trt.age_at_treatment_end 
END AS "Age at Treatment End",
    trt.treatment_type AS "Treatment Type",
    REPLACE(trt.treatment_agent, ';', ', ') AS "Treatment Agent",
    std.dbgap_accession AS "dbGaP Accession"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_diagnoses dgn ON prt.id = dgn."participant.id"
LEFT JOIN 
    df_treatments trt ON prt.id = trt."participant.id"
LEFT JOIN 
    df_treatment_resp trr ON prt.id = trr."participant.id"
LEFT JOIN 
    df_survival srv ON prt.id = srv."participant.id"
LEFT JOIN 
    df_reference_files rfs ON std.id = rfs."study.id"
WHERE 
    std.dbgap_accession = 'phs000466' AND dgn.disease_phase = 'Initial Diagnosis'
ORDER BY 
    trt.treatment_id ASC
LIMIT 100;
"@

$ws.Range("B5").Value = $treatmentQuery

# Nudge the font identity on the Diagnosis (B4) and Treatment (B5) query
# cells so they pick up a distinct (but visually identical, 12pt) font
# record/style, matching the re-formatting captured in the source diff.
$ws.Range("B4").Font.ThemeColor = 1
$ws.Range("B5").Font.ThemeColor = 1

# Move the active selection to C5 (also reflected in the saved view state).
$ws.Range("C5").Select()
